$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text would otherwise be auto-parsed as a number by Excel;
# force them to stay as text so the stored cell type/content matches the source data.
$textCells = @(
    "D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", 
    "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", 
    "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", 
    "D45", "D46", "D47", "D48", "D49", "D50", "D51" 
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.182.76"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "1.753.42"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "239.02"
$ws.Range("E5").Value = "  +3.99%  "
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "0.5290"
$ws.Range("E7").Value = "  +2.64%  "
$ws.Range("D8").Value = "0.2827"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "0.06209"
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").Value = "1.747.29"
$ws.Range("E10").Value = "  +0.16%  "
$ws.Range("D11").Value = "0.07217"
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("D12").Value = "15.58"
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").Value = "0.6503"
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("D14").Value = "4.648"
$ws.Range("E14").Value = "  +3.07%  "
$ws.Range("D15").Value = "78.91"
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("D16").Value = "0.9996"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "0.9988"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").Value = "26.066.20"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").Value = "11.87"
$ws.Range("E19").Value = "  +3.18%  "
$ws.Range("D20").Value = "0.000006770"
$ws.Range("E20").Value = "  +2.39%  "
$ws.Range("D21").Value = "1.971.58"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "4.357"
$ws.Range("E22").Value = "  +5.87%  "
$ws.Range("D23").Value = "8.776"
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("D24").Value = "5.263"
$ws.Range("E24").Value = "  +2.47%  "
$ws.Range("D25").Value = "139.65"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "1.524"
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("D27").Value = "15.39"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").Value = "105.16"
$ws.Range("E29").Value = "  +2.22%  "
$ws.Range("D30").Value = "0.08340"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").Value = "3.835"
$ws.Range("E31").Value = "  +5.64%  "
$ws.Range("D32").Value = "3.674"
$ws.Range("E32").Value = "  +7.47%  "
$ws.Range("D33").Value = "0.04629"
$ws.Range("E33").Value = "  +5.30%  "
$ws.Range("D34").Value = "2.647"
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("D35").Value = "1.027"
$ws.Range("E35").Value = "  +5.03%  "
$ws.Range("D36").Value = "0.6382"
$ws.Range("E36").Value = "  +5.01%  "
$ws.Range("D37").Value = "2.710"
$ws.Range("E37").Value = "  +1.27%  "
$ws.Range("D38").Value = "0.01631"
$ws.Range("E38").Value = "  +4.15%  "
$ws.Range("D39").Value = "1.991"
$ws.Range("E39").Value = "  +3.25%  "
$ws.Range("D40").Value = "0.9991"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "102.54"
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").Value = "0.3965"
$ws.Range("E42").Value = "  +3.21%  "
$ws.Range("D43").Value = "0.7556"
$ws.Range("E43").Value = "  +4.11%  "
$ws.Range("D44").Value = "5.076"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").Value = "0.1159"
$ws.Range("E45").Value = "  +4.60%  "
$ws.Range("D46").Value = "6.442"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("D47").Value = "0.05359"
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("D48").Value = "31.25"
$ws.Range("D49").Value = "54.70"
$ws.Range("E49").Value = "  +4.03%  "
$ws.Range("D50").Value = "0.3496"
$ws.Range("E50").Value = "  +2.93%  "
$ws.Range("D51").Value = "7.612"
$ws.Range("E51").Value = "  +1.34%  "
